$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add footnote/help text cells next to the "role" entry (row 4) and the
# "Organization_name" entry under the Visiting/Co-supervised Students group
# (row 12). These land as new shared strings and plain (unstyled) cells.
$ws.Range("C4").Value = "** If you are a student use : [degree] student, institution "
$ws.Range("C12").Value = "** administrative category, for after receipt."

# Leave the active selection on the newly annotated cell, matching the
# saved view state of the edited workbook.
$ws.Range("C12").Select()
